$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
  @{Row=2; D="57.349.54"; E="  +0.38%  "},
  @{Row=3; D="2.364.11"; E="  +1.54%  "},
  @{Row=4; D="1.00"; E="  -0.19%  "},
  @{Row=5; D="520.19"; E="  +0.24%  "},
  @{Row=6; D="136.00"; E="  +1.04%  "},
  @{Row=7; D=$null; E="  -0.33%  "},
  @{Row=8; D=$null; E="  +0.29%  "},
  @{Row=9; D=$null; E="  -0.70%  "},
  @{Row=10; D="5.49"; E="  +4.49%  "},
  @{Row=11; D=$null; E="  -0.89%  "},
  @{Row=12; D=$null; E="  -0.02%  "},
  @{Row=13; D="24.40"; E="  +2.04%  "},
  @{Row=14; D="2.786.53"; E="  +0.36%  "},
  @{Row=15; D="57.337.06"; E="  +0.62%  "},
  @{Row=16; D="0.0000136"; E="  +0.05%  "},
  @{Row=17; D="2.377.51"; E="  +0.92%  "},
  @{Row=18; D="10.58"; E="  +0.30%  "},
  @{Row=19; D="329.95"; E="  +2.09%  "},
  @{Row=20; D=$null; E="  -1.26%  "},
  @{Row=21; D="6.70"; E="  -0.07%  "},
  @{Row=22; D=$null; E="  -0.04%  "},
  @{Row=23; D="61.45"; E="  +0.13%  "},
  @{Row=24; D="8.94"; E="  +14.66%  "},
  @{Row=25; D=$null; E="  +2.99%  "},
  @{Row=26; D="0.994"; E="  -0.33%  "},
  @{Row=27; D="1.34"; E="  +10.93%  "},
  @{Row=28; D="0.0₃0749"; E="  +1.12%  "},
  @{Row=29; D=$null; E="  +1.32%  "},
  @{Row=30; D="166.42"; E="  -3.25%  "},
  @{Row=31; D="6.27"; E="  -0.46%  "},
  @{Row=32; D="18.61"; E="  +1.19%  "},
  @{Row=33; D=$null; E="  -0.05%  "},
  @{Row=34; D="1.30"; E="  +3.36%  "},
  @{Row=35; D="0.994"; E="  -0.23%  "},
  @{Row=36; D="0.919"; E="  -3.99%  "},
  @{Row=37; D=$null; E="  +0.25%  "},
  @{Row=38; D=$null; E="  +6.10%  "},
  @{Row=39; D="38.86"; E="  +3.42%  "},
  @{Row=40; D="149.74"; E="  +6.77%  "},
  @{Row=41; D="0.388"; E="  +1.23%  "},
  @{Row=44; D="5.25"; E="  +1.66%  "},
  @{Row=45; D="0.0939"; E="  +0.72%  "},
  @{Row=46; D="0.0509"; E="  -0.60%  "},
  @{Row=47; D=$null; E="  +0.67%  "},
  @{Row=48; D=$null; E="  +5.24%  "},
  @{Row=49; D="0.0219"; E="  +1.52%  "},
  @{Row=50; D="17.70"; E="  +4.30%  "},
  @{Row=51; D=$null; E="  -5.25%  "}
)

foreach ($item in $changes) {
  if ($item.D -ne $null) {
    $ws.Cells.Item($item.Row, 4).Value = "'" + $item.D
  }
  if ($item.E -ne $null) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
  }
}

# Row 42/43 swap: Bittensor <-> Filecoin (including link URLs, price, volume)
$ws.Cells.Item(42, 2).Value = "Filecoin"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(42, 4).Value = "'3.66"
$ws.Cells.Item(42, 5).Value = "  +1.74%  "
$ws.Cells.Item(43, 2).Value = "Bittensor"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(43, 4).Value = "'290.40"
$ws.Cells.Item(43, 5).Value = "  +4.68%  "

Write-Host "Applied cryptos update"
